# Generate Report for Handback
# Replaces the old source-file UUID based filenames / hashes / timestamps
# with the new ones produced by a fresh handback run, across the
# Overview, zh-cn and de-de worksheets (cell values + hyperlink display text).

$wb = $excel.ActiveWorkbook

$oldMdA = "d7b4059c-1722-4798-b508-eab6fed8d24f.md"
$oldMdB = "f52a5959-f47d-4956-a8d3-e78bc9049c42.md"
$newMdA = "92734240-b06a-4edf-80fc-7c941305c682.md"
$newMdB = "ffffe431486b-1bb4-478b-bfb7-84d9e87df2d1.md"

$newXlfZh = "92734240-b06a-4edf-80fc-7c941305c682.2a6680a7e3c6e662617f3509a381d64b5dd388c2.zh-cn.xlf"
$newXlfDe = "92734240-b06a-4edf-80fc-7c941305c682.2a6680a7e3c6e662617f3509a381d64b5dd388c2.de-de.xlf"

$newHandoffZh = "2016-03-14 03:44:11"
$newHandbackZh = "2016-03-14 03:44:28"
$newHandoffDe = "2016-03-14 03:44:13"
$newHandbackDe = "2016-03-14 03:44:32"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdA
$wsOverview.Range("A3").Value = $newMdB

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ac69bb2ecfdaf7e421f00889bb784d9c0d2ee6b5/e2e/$newMdA", "", "", $newMdA) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ac69bb2ecfdaf7e421f00889bb784d9c0d2ee6b5/e2e/$newMdB", "", "", $newMdB) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMdA
$wsZh.Range("D2").Value = $newXlfZh
$wsZh.Range("E2").Value = $newHandoffZh
$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("H2").Value = $newHandbackZh

$wsZh.Range("A3").Value = $newMdB
$wsZh.Range("D3").Value = $newXlfZh
$wsZh.Range("E3").Value = $newHandoffZh
$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("H3").Value = $newHandbackZh

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ac69bb2ecfdaf7e421f00889bb784d9c0d2ee6b5/e2e/$newMdA", "", "", $newMdA) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/ac69bb2ecfdaf7e421f00889bb784d9c0d2ee6b5/e2e/$newMdA", "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/511f86e20b01f225eaec4434c7bcd27fa792da5e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfZh", "", "", $newXlfZh) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/23652eeb1c0a51a49994ca88674eca79a903ea51/e2e/$newMdA", "", "", $newMdA) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6ad070f08b996b918e52b6f024a2c7cad426ba66/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfZh", "", "", $newXlfZh) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ac69bb2ecfdaf7e421f00889bb784d9c0d2ee6b5/e2e/$newMdB", "", "", $newMdB) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/ac69bb2ecfdaf7e421f00889bb784d9c0d2ee6b5/e2e/$newMdB", "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/511f86e20b01f225eaec4434c7bcd27fa792da5e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfZh", "", "", $newXlfZh) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/23652eeb1c0a51a49994ca88674eca79a903ea51/e2e/$newMdB", "", "", $newMdB) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6ad070f08b996b918e52b6f024a2c7cad426ba66/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfZh", "", "", $newXlfZh) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMdA
$wsDe.Range("D2").Value = $newXlfDe
$wsDe.Range("E2").Value = $newHandoffDe
$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("H2").Value = $newHandbackDe

$wsDe.Range("A3").Value = $newMdB
$wsDe.Range("D3").Value = $newXlfDe
$wsDe.Range("E3").Value = $newHandoffDe
$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("H3").Value = $newHandbackDe

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ac69bb2ecfdaf7e421f00889bb784d9c0d2ee6b5/e2e/$newMdA", "", "", $newMdA) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/ac69bb2ecfdaf7e421f00889bb784d9c0d2ee6b5/e2e/$newMdA", "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43e3c36be1e1506ed736e8e6e17495aea645f599/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDe", "", "", $newXlfDe) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5c7e3bd7efa82d5c8431d912b034cd549aa3627b/e2e/$newMdA", "", "", $newMdA) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d9d957e91b19a6847edf0beaae9aac84adce75f2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDe", "", "", $newXlfDe) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ac69bb2ecfdaf7e421f00889bb784d9c0d2ee6b5/e2e/$newMdB", "", "", $newMdB) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/ac69bb2ecfdaf7e421f00889bb784d9c0d2ee6b5/e2e/$newMdB", "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43e3c36be1e1506ed736e8e6e17495aea645f599/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDe", "", "", $newXlfDe) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5c7e3bd7efa82d5c8431d912b034cd549aa3627b/e2e/$newMdB", "", "", $newMdB) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d9d957e91b19a6847edf0beaae9aac84adce75f2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDe", "", "", $newXlfDe) | Out-Null
